# Rename headers on existing sheets
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, centered, bordered) from an existing header cell
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$arr = New-Object 'object[,]' 50,4
$arr[0,0]=45207.99999999999; $arr[0,1]=41; $arr[0,2]=-90.66872078194; $arr[0,3]=174.891566071575
$arr[1,0]=45214.99999999999; $arr[1,1]=43; $arr[1,2]=-93.99846192844886; $arr[1,3]=168.0701914670582
$arr[2,0]=45221.99999999999; $arr[2,1]=44; $arr[2,2]=-83.57839182582784; $arr[2,3]=176.8646399645672
$arr[3,0]=45228.99999999999; $arr[3,1]=45; $arr[3,2]=-91.85006568799612; $arr[3,3]=178.0128804777577
$arr[4,0]=45235.99999999999; $arr[4,1]=46; $arr[4,2]=-90.7895138738582; $arr[4,3]=176.3675736060794
$arr[5,0]=45242.99999999999; $arr[5,1]=47; $arr[5,2]=-87.55993462089722; $arr[5,3]=177.0847803628496
$arr[6,0]=45249.99999999999; $arr[6,1]=48; $arr[6,2]=-91.36444145845201; $arr[6,3]=170.4736356690826
$arr[7,0]=45256.99999999999; $arr[7,1]=49; $arr[7,2]=-81.82807524176926; $arr[7,3]=177.3277681823568
$arr[8,0]=45263.99999999999; $arr[8,1]=50; $arr[8,2]=-74.8443843618642; $arr[8,3]=180.5896345634308
$arr[9,0]=45326.99999999999; $arr[9,1]=59; $arr[9,2]=-71.51203990305179; $arr[9,3]=199.6889375421677
$arr[10,0]=45333.99999999999; $arr[10,1]=60; $arr[10,2]=-75.22850737013344; $arr[10,3]=193.3903562888575
$arr[11,0]=45340.99999999999; $arr[11,1]=61; $arr[11,2]=-73.96353979985133; $arr[11,3]=189.7629362819077
$arr[12,0]=45347.99999999999; $arr[12,1]=62; $arr[12,2]=-53.35898400294813; $arr[12,3]=195.1846969785103
$arr[13,0]=45354.99999999999; $arr[13,1]=63; $arr[13,2]=-73.70797928640624; $arr[13,3]=195.0876879178187
$arr[14,0]=45361.99999999999; $arr[14,1]=64; $arr[14,2]=-71.44544026354656; $arr[14,3]=199.1708080466981
$arr[15,0]=45368.99999999999; $arr[15,1]=65; $arr[15,2]=-67.44402614335915; $arr[15,3]=188.7561425033183
$arr[16,0]=45375.99999999999; $arr[16,1]=66; $arr[16,2]=-63.71028034757629; $arr[16,3]=210.477571810354
$arr[17,0]=45382.99999999999; $arr[17,1]=67; $arr[17,2]=-66.56509509918267; $arr[17,3]=193.1316571557421
$arr[18,0]=45389.99999999999; $arr[18,1]=68; $arr[18,2]=-59.85306365204585; $arr[18,3]=199.8818085645502
$arr[19,0]=45396.99999999999; $arr[19,1]=69; $arr[19,2]=-58.1732889124111; $arr[19,3]=197.8064231544472
$arr[20,0]=45417.99999999999; $arr[20,1]=73; $arr[20,2]=-53.32978324617025; $arr[20,3]=208.4666913118727
$arr[21,0]=45424.99999999999; $arr[21,1]=74; $arr[21,2]=-56.80165420090386; $arr[21,3]=209.2611866538984
$arr[22,0]=45431.99999999999; $arr[22,1]=75; $arr[22,2]=-52.67420089854738; $arr[22,3]=215.4474330129594
$arr[23,0]=45438.99999999999; $arr[23,1]=76; $arr[23,2]=-57.58201748027749; $arr[23,3]=196.8076068461238
$arr[24,0]=45445.99999999999; $arr[24,1]=77; $arr[24,2]=-50.45890201387516; $arr[24,3]=215.5618438270902
$arr[25,0]=45452.99999999999; $arr[25,1]=78; $arr[25,2]=-52.47450584385074; $arr[25,3]=200.4296319525432
$arr[26,0]=45459.99999999999; $arr[26,1]=79; $arr[26,2]=-51.475927093919; $arr[26,3]=210.5351557628362
$arr[27,0]=45466.99999999999; $arr[27,1]=80; $arr[27,2]=-50.40850395582372; $arr[27,3]=215.5766864360512
$arr[28,0]=45473.99999999999; $arr[28,1]=81; $arr[28,2]=-49.41608788040549; $arr[28,3]=203.5203032106204
$arr[29,0]=45480.99999999999; $arr[29,1]=82; $arr[29,2]=-42.67636699510993; $arr[29,3]=217.6603359109833
$arr[30,0]=45487.99999999999; $arr[30,1]=83; $arr[30,2]=-45.06397922377502; $arr[30,3]=209.5252557761727
$arr[31,0]=45501.99999999999; $arr[31,1]=85; $arr[31,2]=-45.44182068346369; $arr[31,3]=213.8148545375986
$arr[32,0]=45508.99999999999; $arr[32,1]=86; $arr[32,2]=-40.73358664183184; $arr[32,3]=215.2221277711008
$arr[33,0]=45515.99999999999; $arr[33,1]=87; $arr[33,2]=-45.72822009436253; $arr[33,3]=216.0790988886101
$arr[34,0]=45522.99999999999; $arr[34,1]=88; $arr[34,2]=-39.35325942671523; $arr[34,3]=221.345164801174
$arr[35,0]=45536.99999999999; $arr[35,1]=90; $arr[35,2]=-36.99324315410711; $arr[35,3]=221.6924319656547
$arr[36,0]=45543.99999999999; $arr[36,1]=91; $arr[36,2]=-36.06313040680885; $arr[36,3]=234.4421399100779
$arr[37,0]=45557.99999999999; $arr[37,1]=93; $arr[37,2]=-42.27383925292838; $arr[37,3]=231.5665718019249
$arr[38,0]=45578.99999999999; $arr[38,1]=96; $arr[38,2]=-33.3422965815422; $arr[38,3]=224.2705643306762
$arr[39,0]=45585.99999999999; $arr[39,1]=97; $arr[39,2]=-33.84212308374682; $arr[39,3]=233.7697336881544
$arr[40,0]=45592.99999999999; $arr[40,1]=99; $arr[40,2]=-47.19515079387457; $arr[40,3]=217.1758518870413
$arr[41,0]=45613.99999999999; $arr[41,1]=102; $arr[41,2]=-21.63465718374052; $arr[41,3]=239.4094998348182
$arr[42,0]=45620.99999999999; $arr[42,1]=103; $arr[42,2]=-21.0357489690597; $arr[42,3]=235.3049247889845
$arr[43,0]=45627.99999999999; $arr[43,1]=104; $arr[43,2]=-31.09233637739603; $arr[43,3]=235.4587894895422
$arr[44,0]=45634.99999999999; $arr[44,1]=105; $arr[44,2]=-22.45666856062826; $arr[44,3]=232.2254843558869
$arr[45,0]=45641.99999999999; $arr[45,1]=106; $arr[45,2]=-31.59104618825571; $arr[45,3]=233.1526296364119
$arr[46,0]=45648.99999999999; $arr[46,1]=107; $arr[46,2]=-27.21679359472976; $arr[46,3]=229.5181224550217
$arr[47,0]=45655.99999999999; $arr[47,1]=108; $arr[47,2]=-23.86811097302571; $arr[47,3]=250.0158142063807
$arr[48,0]=45662.99999999999; $arr[48,1]=109; $arr[48,2]=-19.64529592492675; $arr[48,3]=242.0598031075537
$arr[49,0]=45669.99999999999; $arr[49,1]=110; $arr[49,2]=-17.79038444052197; $arr[49,3]=242.5510016323596

$dataRange = $wsForecast.Range("A2:D51")
$dataRange.Value = $arr

# Copy date-format style from an existing date cell onto the new ds column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A51").PasteSpecial(-4122)

# Restore the originally active sheet/selection
[void]$wsWeekly.Activate()
[void]$wsWeekly.Range("A1").Select()
